$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Sculptor worked on a statue of MC" -> "Sculptor secretly worked on a
#    statue of MC", split across three runs: "Sculptor" / " secretly" /
#    " worked on a statue of MC".
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Sculptor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(0)
$ins1 = $r1.Duplicate
$ins1.InsertAfter(" secretly")
# Force a run break around the freshly inserted text (toggling a format on
# and back off keeps it a distinct <w:r> even though the effective
# formatting ends up identical to its neighbours).
$ins1.Font.Bold = $true
$ins1.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2) "MC can't face a portrayal of himself, pushed the statue onto the
#    ground" -> "...onto the floor", split into the unchanged lead-in run
#    plus a new "floor" run.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("MC can" + [char]8217 + "t face a portrayal of himself, pushed the statue onto the ground", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2b = $r2.Duplicate
$r2b.Find.Execute("ground", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2b.Text = "floor"
$r2b.Font.Bold = $true
$r2b.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) "The statue is destroyed and the fragments lay on the ground" ->
#    "The statue is destroyed, shattered into dozens of fragments", split
#    into three runs: "The statue is destroyed" / ", shattered into dozens
#    of f" / "ragments". The trailing bookmark ("_GoBack") that used to sit
#    at the end of this paragraph's text moves to the end of the final
#    (all-tabs) paragraph instead.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("The statue is destroyed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail = $d.Content
$tail.Find.Execute(" and the fragments lay on the ground", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Text = ", shattered into dozens of fragments"

$split1 = $d.Content
$split1.Find.Execute(", shattered into dozens of f", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$split1.Font.Bold = $true
$split1.Font.Bold = $false

$split2 = $d.Content
$split2.Find.Execute("ragments", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$split2.Font.Bold = $true
$split2.Font.Bold = $false

# ---------------------------------------------------------------------------
# Move the "_GoBack" bookmark from the end of the "statue is destroyed"
# paragraph to the end of the very last (all-tabs) paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPos = $lastPara.Range.Duplicate
$endPos.Collapse(0)
$endPos.MoveEnd(1, -1) | Out-Null
$endPos.Collapse(0)

# Creating a bookmark collapsed right at the very end of the document's
# content is unreliable, so temporarily extend the document, anchor the
# bookmark, then trim the placeholder back off again.
$endPos.InsertAfter("X")
$endPos.Collapse(0)
$endPos.MoveStart(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $endPos) | Out-Null
$placeholder = $d.Range($endPos.Start, $endPos.Start + 1)
$placeholder.Delete() | Out-Null
